$d = $word.ActiveDocument

# --- Change 1: merge "Católico/Protestante: " + "Religion;" runs into one run ---
$d.Content.Find.Execute("Católico/Protestante: Religion;", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Católico/Protestante: Religion;", 2)

# --- Change 2: insert the new "INTERFACE" block right after the "Ações (TODO):"
#     paragraph, before the "Unique name ..." list paragraph ---
$findRng = $d.Content
$found = $findRng.Find.Execute("(TODO):", $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)

if ($found) {
    $anchorPara = $findRng.Paragraphs(1)
    $anchorStart = $anchorPara.Range.Start

    # Resolve the numeric Paragraphs() index for the anchor paragraph so we can
    # keep stepping forward through freshly inserted paragraphs by index.
    $anchorIndex = 0
    $count = 0
    foreach ($p in $d.Paragraphs) {
        $count = $count + 1
        if ($p.Range.Start -eq $anchorStart) {
            $anchorIndex = $count
        }
    }

    if ($anchorIndex -gt 0) {
        $lines = @(
            "INTERFACE",
            "- desaparecer a tela inicial de adicionar jogadores",
            "- mostrar as cartas de cada jogador na mesa",
            "- indicar o jogador atual",
            "- mostrar ações possíveis para o jogador atual"
        )

        $insertIndex = $anchorIndex
        foreach ($line in $lines) {
            $anchor = $d.Paragraphs($insertIndex)
            $rng = $anchor.Range
            $rng.Collapse(0)
            $rng.InsertParagraphAfter()
            $insertIndex = $insertIndex + 1
            $newPara = $d.Paragraphs($insertIndex)
            $newPara.Range.InsertAfter($line)
        }

        # Trailing empty bold paragraph - insert then strip the placeholder
        # character back out so no <w:r> survives in the final XML (matching
        # the blank bold paragraphs used elsewhere in this document).
        $anchor = $d.Paragraphs($insertIndex)
        $rng = $anchor.Range
        $rng.Collapse(0)
        $rng.InsertParagraphAfter()
        $insertIndex = $insertIndex + 1
        $emptyPara = $d.Paragraphs($insertIndex)
        $emptyPara.Range.InsertAfter("X")
        $cleanupRng = $d.Range($emptyPara.Range.Start, $emptyPara.Range.Start + 1)
        $cleanupRng.Delete()
    }
}

# --- Change 3: merge "Unique name – dois " + "jogadores não podem ter o mesmo nome;" runs ---
$d.Content.Find.Execute("Unique name – dois jogadores não podem ter o mesmo nome;", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Unique name – dois jogadores não podem ter o mesmo nome;", 2)

# --- Change 4: merge "V" + "erificar regulamentação;" runs ---
$d.Content.Find.Execute("Verificar regulamentação;", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Verificar regulamentação;", 2)
